# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 24; I = "sv"; J = "Statement-opinion" },
    @{ Row = 25; I = "sv"; J = "Statement-opinion" },
    @{ Row = 26; I = "aa"; J = "Agree/Accept" },
    @{ Row = 27; I = "aa"; J = "Agree/Accept" },
    @{ Row = 31; I = "aa"; J = "Agree/Accept" },
    @{ Row = 36; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 45; I = "aa"; J = "Agree/Accept" },
    @{ Row = 46; I = "aa"; J = "Agree/Accept" },
    @{ Row = 48; I = "sv"; J = "Statement-opinion" },
    @{ Row = 49; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 53; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 55; I = "aa"; J = "Agree/Accept" },
    @{ Row = 56; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 57; I = "aa"; J = "Agree/Accept" },
    @{ Row = 58; I = "aa"; J = "Agree/Accept" },
    @{ Row = 79; I = "ba"; J = "Appreciation" },
    @{ Row = 80; I = "aa"; J = "Agree/Accept" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
